$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97:128 down to 98:129.
$ws.Rows.Item(97).Insert()

# Populate the new row 97 with the weekly record (same template as the
# surrounding rows, new date + volume/price figures).
$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 45119
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = 100112005
$ws.Range("G97").Value = "Puerro"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 70
$ws.Range("K97").Value = 8000
$ws.Range("L97").Value = 8000
$ws.Range("M97").Value = 8000
$ws.Range("N97").Value = "$/paquete 20 unidades"
$ws.Range("O97").Value = "Provincia de Chacabuco"
$ws.Range("P97").Value = 400
$ws.Range("Q97").Value = 20
$ws.Range("R97").Value = "Hortaliza"
